$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header (H1, bold/centered/bordered)
# onto the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data for rows 2-11: column I is always 1, column J mirrors column H (IP)
for ($row = 2; $row -le 11; $row++) {
    $hVal = [double]$ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hVal
}
